{"js": "// CU-24 Cancelar merma \u2014 apply the described corrections.\n\nconst body = context.document.body;\n\n// 1) Trazabilidad: \"FRQ- \" -> \"FRQ-30 \"\nconst frqResults = body.search(\"FRQ- \", { matchCase: true });\nfrqResults.load(\"items\");\nawait context.sync();\nif (frqResults.items.length > 0) {\n  const afterDash = frqResults.items[0].insertText(\"FRQ-\", \"Replace\");\n  await context.sync();\n  const afterNum = afterDash.insertText(\"30\", \"After\");\n  await context.sync();\n  afterNum.insertText(\" \", \"After\");\n  await context.sync();\n}\n\n// 2) Actor(es): \"Administrador del supermercado, \" -> \"Administrador, \"\nconst actorResults = body.search(\"Administrador del supermercado, \", { matchCase: true });\nactorResults.load(\"items\");\nawait context.sync();\nif (actorResults.items.length > 0) {\n  actorResults.items[0].insertText(\"Administrador, \", \"Replace\");\n  await context.sync();\n}\n\n// 3) Disparador: reword the trigger sentence.\nconst selResults = body.search(\"selecciona de lista la MERMA que desea eliminar y \", { matchCase: true });\nselResults.load(\"items\");\nawait context.sync();\nif (selResults.items.length > 0) {\n  selResults.items[0].insertText(\"selecciona de lista la MERMA que desea cancelar y \", \"Replace\");\n  await context.sync();\n}\n\nconst clickResults = body.search(\"hace clic en el bot\u00f3n \\u201cCancelar merma\\u201d\", { matchCase: true });\nclickResults.load(\"items\");\nawait context.sync();\nif (clickResults.items.length > 0) {\n  clickResults.items[0].insertText(\"hace clic en el bot\u00f3n \\u201cCancelar\\u201d\", \"Replace\");\n  await context.sync();\n}\n\n// 4) Add a new \"Extiende\" row at the end of the table.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.addRows(\"End\", 1, [[\"Extiende\", \"CU-38 Consultar merma.\"]]);\n  await context.sync();\n}\n", "ps1": "# CU-24 Cancelar merma \u2014 apply the described corrections.\n$d = $word.ActiveDocument\n\n$LDQ = [char]0x201C   # \u201c\n$RDQ = [char]0x201D   # \u201d\n$OACUTE = [char]0x00F3  # \u00f3\n\n# 1) Trazabilidad: \"FRQ- \" -> \"FRQ-30 \"\n$find1 = $d.Content.Find\n$find1.Execute(\"FRQ- \", $false, $false, $false, $false, $false, $true, 1, $false, \"FRQ-30 \", 2) | Out-Null\n\n# 2) Actor(es): \"Administrador del supermercado, \" -> \"Administrador, \"\n$find2 = $d.Content.Find\n$find2.Execute(\"Administrador del supermercado, \", $false, $false, $false, $false, $false, $true, 1, $false, \"Administrador, \", 2) | Out-Null\n\n# 3) Disparador: reword the trigger sentence.\n$find3 = $d.Content.Find\n$find3.Execute(\"selecciona de lista la MERMA que desea eliminar y \", $false, $false, $false, $false, $false, $true, 1, $false, \"selecciona de lista la MERMA que desea cancelar y \", 2) | Out-Null\n\n$boton = \"bot\" + $OACUTE + \"n\"\n$find4 = $d.Content.Find\n$find4.Execute(\"hace clic en el $boton \" + $LDQ + \"Cancelar merma\" + $RDQ, $false, $false, $false, $false, $false, $true, 1, $false, \"hace clic en el $boton \" + $LDQ + \"Cancelar\" + $RDQ, 2) | Out-Null\n\n# 4) Add a new \"Extiende\" row at the end of the table.\n$t = $d.Tables.Item(1)\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"Extiende\"\n$newRow.Cells.Item(2).Range.Text = \"CU-38 Consultar merma.\"\n"}
